$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns(1)

function Set-NameBySku($sku, $newName) {
    $cell = $col.Find($sku)
    if ($cell -ne $null) {
        $ws.Cells.Item($cell.Row, 2).Value = $newName
    }
}

Set-NameBySku "SCF-100.PRO8340"        "Astonish 6in1 ASTONISH BODY & SOUL MINI SHOWER GEL GIFT PACK"
Set-NameBySku "SCF-SCRC-BASKIN.KING"   "BaskinRobbin Single King Scoop"
Set-NameBySku "SCF-SCRC-BASKIN.R"      "BaskinRobbin Single Regular Scoop"
Set-NameBySku "SCF-100.ETX9807"        "Cento 2 Way Heavy Cable Reel 7m"
Set-NameBySku "SCF-100.ECOGS231REC"    "Cento Ecocen GLASS STORAGE"
Set-NameBySku "SCF-100.TR606P35"       "FELIZTRIP TRAVEL ADAPTOR. 2x USB C and 2X USB A Fast Charging PD 35W"
Set-NameBySku "SCF-100.GLADELEMON_T/P_MD" "GLADE Scent Gel TwinPack Lemon"
Set-NameBySku "SCF-100.GLO800ML_LEMON" "GLO 800ML Lemon"
Set-NameBySku "SCF-100.WB55"           "JIMMY DUST MITE KILLER"
Set-NameBySku "SCF-100.VACUUMSTAND"    "JIMMY VACUUM STAND"
Set-NameBySku "SCF-100.BLC129"         "Khind 500W FOOD PROCESSOR"
Set-NameBySku "SCF-100.CF683DC"        "Khind CEILING FAN"
Set-NameBySku "SCF-100.FPC900"         "Khind CHOPPER"
Set-NameBySku "SCF-100.MRMUSLE_CITRUS" "Mr.Muscle TOILET BOWL CLEANER - CITRUS"
Set-NameBySku "SCF-100.MRMUSLE_MAR"    "Mr.Muscle TOILET BOWL CLEANER - MARINE"
Set-NameBySku "SCF-100.MRMUSLE_PINE"   "Mr.Muscle TOILET BOWL CLEANER - PINE"
Set-NameBySku "SCF-100.HD3064/62"      "Philips 0.54L RICE COOKER"
Set-NameBySku "SCF-100.T3"             "RussellTaylors Digital Bread Toaster"
Set-NameBySku "SCF-100.SAFE95587"      "PanzerGlass SAFE Apple iPhone 2023 6.1`" Pro UWF"
Set-NameBySku "SCF-100.SAFE95589"      "PanzerGlass SAFE Apple iPhone 2023 6.7`" Pro UWF"

$ws.Range("B131").Select()
